$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell value as plain text, preserving the cells original
# (default) style -- Range.Value auto-coerces numeric-looking strings into
# real numbers (e.g. "1.00" -> 1), which loses formatting the source data
# relies on. Forcing a Text number format for the assignment keeps the
# literal string, then resetting the Style back to "Normal" removes the
# now-unneeded explicit text-format style so the cells style index is
# left exactly as it was before (no stray "s=" attribute in the XML).
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range('D2') '28.050.85'
Set-TextValue $ws.Range('E2') '  +3.27%  '
Set-TextValue $ws.Range('D3') '1.688.38'
Set-TextValue $ws.Range('E3') '  +0.42%  '
Set-TextValue $ws.Range('D4') '1.00'
Set-TextValue $ws.Range('E4') '  -0.31%  '
Set-TextValue $ws.Range('D5') '216.67'
Set-TextValue $ws.Range('E5') '  +0.82%  '
Set-TextValue $ws.Range('D6') '0.520'
Set-TextValue $ws.Range('E6') '  +0.41%  '
Set-TextValue $ws.Range('D7') '1.00'
Set-TextValue $ws.Range('E7') '  -0.29%  '
Set-TextValue $ws.Range('E8') '  +6.68%  '
Set-TextValue $ws.Range('D9') '0.265'
Set-TextValue $ws.Range('E9') '  +1.89%  '
Set-TextValue $ws.Range('E10') '  +0.40%  '
Set-TextValue $ws.Range('D11') '0.0885'
Set-TextValue $ws.Range('E11') '  -0.60%  '
Set-TextValue $ws.Range('D12') '1.927.06'
Set-TextValue $ws.Range('E12') '  +0.36%  '
Set-TextValue $ws.Range('D13') '1.689.27'
Set-TextValue $ws.Range('E13') '  +0.46%  '
Set-TextValue $ws.Range('E14') '  -0.08%  '
Set-TextValue $ws.Range('E15') '  +0.50%  '
Set-TextValue $ws.Range('D16') '66.89'
Set-TextValue $ws.Range('E16') '  +0.24%  '
Set-TextValue $ws.Range('D17') '250.46'
Set-TextValue $ws.Range('E17') '  +6.25%  '
Set-TextValue $ws.Range('D18') '28.008.30'
Set-TextValue $ws.Range('E18') '  +3.14%  '
Set-TextValue $ws.Range('D19') '0.0₃0742'
Set-TextValue $ws.Range('E19') '  +0.48%  '
Set-TextValue $ws.Range('E20') '  -2.36%  '
Set-TextValue $ws.Range('E21') '  -0.15%  '
Set-TextValue $ws.Range('E22') '  -0.13%  '
Set-TextValue $ws.Range('D23') '9.57'
Set-TextValue $ws.Range('E23') '  +0.42%  '
Set-TextValue $ws.Range('E24') '  -1.61%  '
Set-TextValue $ws.Range('D25') '147.61'
Set-TextValue $ws.Range('E25') '  +0.53%  '
Set-TextValue $ws.Range('D26') '7.36'
Set-TextValue $ws.Range('E26') '  -0.72%  '
Set-TextValue $ws.Range('E27') '  +0.97%  '
Set-TextValue $ws.Range('E28') '  +0.57%  '
Set-TextValue $ws.Range('E29') '  -0.27%  '
Set-TextValue $ws.Range('E30') '  +6.81%  '
Set-TextValue $ws.Range('D31') '0.0503'
Set-TextValue $ws.Range('E31') '  -0.04%  '
Set-TextValue $ws.Range('E32') '  +0.53%  '
Set-TextValue $ws.Range('D33') '3.19'
Set-TextValue $ws.Range('E33') '  -1.78%  '
Set-TextValue $ws.Range('D34') '1.426.86'
Set-TextValue $ws.Range('E34') '  -7.59%  '
Set-TextValue $ws.Range('E35') '  -2.54%  '
Set-TextValue $ws.Range('D36') '0.945'
Set-TextValue $ws.Range('E36') '  -0.05%  '
Set-TextValue $ws.Range('D37') '2.39'
Set-TextValue $ws.Range('E37') '  -0.02%  '
Set-TextValue $ws.Range('E38') '  -2.07%  '
Set-TextValue $ws.Range('D39') '0.0173'
Set-TextValue $ws.Range('E39') '  +0.60%  '
Set-TextValue $ws.Range('E40') '  -3.22%  '
Set-TextValue $ws.Range('D41') '69.58'
Set-TextValue $ws.Range('E41') '  +0.73%  '
Set-TextValue $ws.Range('D42') '0.999'
Set-TextValue $ws.Range('E42') '  -0.32%  '
Set-TextValue $ws.Range('D43') '5.52'
Set-TextValue $ws.Range('E43') '  -4.41%  '
Set-TextValue $ws.Range('D44') '1.834.62'
Set-TextValue $ws.Range('E44') '  +0.34%  '
Set-TextValue $ws.Range('E45') '  -0.79%  '
Set-TextValue $ws.Range('D46') '0.798'
Set-TextValue $ws.Range('E46') '  +0.93%  '
Set-TextValue $ws.Range('E47') '  +6.01%  '
Set-TextValue $ws.Range('D48') '89.42'
Set-TextValue $ws.Range('E48') '  -0.66%  '
Set-TextValue $ws.Range('D49') '0.0₆0112'
Set-TextValue $ws.Range('E49') '  -0.14%  '
Set-TextValue $ws.Range('E50') '  -0.62%  '
Set-TextValue $ws.Range('D51') '7.87'
Set-TextValue $ws.Range('E51') '  -3.50%  '
